$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1788
$ws.Range("I18").Value = 990
$ws.Range("J18").Value = 3783
$ws.Range("K18").Value = 990
$ws.Range("L18").Value = 3783
$ws.Range("M18").Value = -706
$ws.Range("N18").Value = -4351
$ws.Range("H21").Value = 29502.834
$ws.Range("I21").Value = 29502.834
$ws.Range("K21").Value = 29502.834
$ws.Range("M21").Value = -29034.834
$ws.Range("H23").Value = 29502.834
$ws.Range("I23").Value = 29502.834
$ws.Range("K23").Value = 29502.834
$ws.Range("M23").Value = -29268.834
$ws.Range("H80").Value = 1714.4286
$ws.Range("I80").Value = 1513
$ws.Range("K80").Value = 4539
$ws.Range("M80").Value = -3541
$ws.Range("H83").Value = 1714.4286
$ws.Range("I83").Value = 1513
$ws.Range("K83").Value = 13617
$ws.Range("M83").Value = -8625
$ws.Range("H131").Value = 7499.75
$ws.Range("I131").Value = 3999.5
$ws.Range("J131").Value = 11000
$ws.Range("K131").Value = 11998.5
$ws.Range("L131").Value = 33000
$ws.Range("M131").Value = -6958.5
$ws.Range("N131").Value = -43080
$ws.Range("H132").Value = 13042.333
$ws.Range("I132").Value = 13042.333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 39126.999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -36596.999
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 2149.0667
$ws.Range("I138").Value = 1176.5714
$ws.Range("K138").Value = 3529.7142
$ws.Range("M138").Value = 1610.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 3104.0833
$ws.Range("I74").Value = 1906.125
$ws.Range("K74").Value = 1906.125
$ws.Range("M74").Value = -1032.125
$ws.Range("H77").Value = 3104.0833
$ws.Range("I77").Value = 1906.125
$ws.Range("K77").Value = 9530.625
$ws.Range("M77").Value = -5162.625
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 50156.875
$ws.Range("I7").Value = 57252.145
$ws.Range("J7").Value = 490
$ws.Range("K7").Value = 57252.145
$ws.Range("L7").Value = 490
$ws.Range("M7").Value = -57139.145
$ws.Range("N7").Value = -716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 183
$ws.Range("I2").Value = 24
$ws.Range("J2").Value = 262.5
$ws.Range("K2").Value = 24
$ws.Range("L2").Value = 262.5
$ws.Range("M2").Value = 89
$ws.Range("N2").Value = -488.5
$ws.Range("H31").Value = 3399
$ws.Range("J31").Value = 3916.5715
$ws.Range("L31").Value = 3916.5715
$ws.Range("N31").Value = -4506.5715
$ws.Range("H34").Value = 3399
$ws.Range("J34").Value = 3916.5715
$ws.Range("L34").Value = 3916.5715
$ws.Range("N34").Value = -4320.5715
$ws.Range("H58").Value = 2078.8333
$ws.Range("I58").Value = 993.5
$ws.Range("K58").Value = 993.5
$ws.Range("M58").Value = -790.5
$ws.Range("H86").Value = 4753.5
$ws.Range("I86").Value = 4753.5
$ws.Range("K86").Value = 4753.5
$ws.Range("M86").Value = -3630.5
$ws.Range("H88").Value = 7000
$ws.Range("J88").Value = 7000
$ws.Range("L88").Value = 7000
$ws.Range("N88").Value = -7812
$ws.Range("H89").Value = 4753.5
$ws.Range("I89").Value = 4753.5
$ws.Range("K89").Value = 23767.5
$ws.Range("M89").Value = -18151.5
$ws.Range("H91").Value = 7000
$ws.Range("J91").Value = 7000
$ws.Range("L91").Value = 7000
$ws.Range("N91").Value = -9808
$ws.Range("H107").Value = 958.5
$ws.Range("J107").Value = 924.3333
$ws.Range("L107").Value = 924.3333
$ws.Range("N107").Value = -4764.3333
$ws.Range("H122").Value = 499.625
$ws.Range("I122").Value = 499.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1498.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 951.125
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2816.6365
$ws.Range("I132").Value = 2779.75
$ws.Range("K132").Value = 8339.25
$ws.Range("M132").Value = -5809.25
$ws.Range("H136").Value = 2078.8333
$ws.Range("I136").Value = 993.5
$ws.Range("K136").Value = 2980.5
$ws.Range("M136").Value = -430.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 250
$ws.Range("I17").Value = 250
$ws.Range("K17").Value = 750
$ws.Range("M17").Value = -581
$ws.Range("H21").Value = 98
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 98
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 294
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -640
$ws.Range("H39").Value = 2125
$ws.Range("J39").Value = 2214.2856
$ws.Range("L39").Value = 6642.8568
$ws.Range("N39").Value = -7230.8568
$ws.Range("H108").Value = 10000
$ws.Range("I108").Value = 10000
$ws.Range("K108").Value = 30000
$ws.Range("M108").Value = -27120
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H118").Value = 2949.75
$ws.Range("J118").Value = 1899
$ws.Range("L118").Value = 5697
$ws.Range("N118").Value = -8183

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 625019.7
$ws.Range("I20").Value = 625019.7
$ws.Range("K20").Value = 625019.7
$ws.Range("M20").Value = -624774.7
$ws.Range("H24").Value = 10000000
$ws.Range("I24").Value = 10000000
$ws.Range("K24").Value = 10000000
$ws.Range("M24").Value = -9999827
$ws.Range("H80").Value = 21600
$ws.Range("J80").Value = 27633.334
$ws.Range("L80").Value = 27633.334
$ws.Range("N80").Value = -29629.334
$ws.Range("H83").Value = 21600
$ws.Range("J83").Value = 27633.334
$ws.Range("L83").Value = 138166.67
$ws.Range("N83").Value = -148150.67
$ws.Range("H107").Value = 502.5
$ws.Range("J107").Value = 507.2
$ws.Range("L107").Value = 507.2
$ws.Range("N107").Value = -4347.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5497.5
$ws.Range("I7").Value = 5497.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5497.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -5385.5
$ws.Range("N7").ClearContents()
$ws.Range("H55").Value = 685.1
$ws.Range("I55").Value = 328.18182
$ws.Range("J55").Value = 1121.3334
$ws.Range("K55").Value = 328.18182
$ws.Range("L55").Value = 1121.3334
$ws.Range("M55").Value = -155.18182
$ws.Range("N55").Value = -1467.3334
$ws.Range("H93").Value = 855
$ws.Range("I93").Value = 818.75
$ws.Range("K93").Value = 818.75
$ws.Range("M93").Value = 429.25
$ws.Range("H126").Value = 5497.5
$ws.Range("I126").Value = 5497.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16492.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14022.5
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336
$ws.Range("H26").Value = 7300000
$ws.Range("I26").Value = 3000000
$ws.Range("K26").Value = 3000000
$ws.Range("M26").Value = -2999707
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 298.4
$ws.Range("I100").Value = 285.5
$ws.Range("K100").Value = 571
$ws.Range("M100").Value = -30
$ws.Range("H113").Value = 2820.3333
$ws.Range("J113").Value = 3200
$ws.Range("L113").Value = 9600
$ws.Range("N113").Value = -13940
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1472.2
$ws.Range("I132").Value = 1428.5
$ws.Range("K132").Value = 4285.5
$ws.Range("M132").Value = -1755.5
$ws.Range("H136").Value = 5480.5557
$ws.Range("I136").Value = 5729.625
$ws.Range("K136").Value = 17188.875
$ws.Range("M136").Value = -14638.875
